$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "angular frequency"
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"

$ws.Range("B1").Select()
